$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal TEXT (matches source data where the
# "Price"/"Volume(1h)" columns are stored as inline strings, not numbers -
# even values that look like plain decimals, e.g. "572.02"). Forcing the
# NumberFormat to Text ("@") before the assignment stops Excel's COM layer
# from auto-coercing numeric-looking strings into real numbers.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# --- Simple Price (D) / Volume(1h) (E) updates ---
Set-TextValue "D2" '64.498.78'
Set-TextValue "E2" '  -2.40%  '
Set-TextValue "D3" '3.179.46'
Set-TextValue "E3" '  -4.23%  '
Set-TextValue "E4" '  +0.00%  '
Set-TextValue "D5" '572.02'
Set-TextValue "E5" '  -2.47%  '
Set-TextValue "D6" '169.55'
Set-TextValue "E6" '  -6.59%  '
Set-TextValue "E7" '  -6.27%  '
Set-TextValue "E8" '  -0.14%  '
Set-TextValue "D9" '3.189.29'
Set-TextValue "E9" '  -3.81%  '
Set-TextValue "E10" '  -3.85%  '
Set-TextValue "D11" '6.85'
Set-TextValue "E11" '  +0.81%  '
Set-TextValue "E12" '  -2.82%  '
Set-TextValue "D13" '3.738.15'
Set-TextValue "E13" '  -4.07%  '
Set-TextValue "D14" '0.129'
Set-TextValue "E14" '  -1.70%  '
Set-TextValue "D15" '64.519.49'
Set-TextValue "E15" '  -2.46%  '
Set-TextValue "D16" '25.33'
Set-TextValue "E16" '  -3.28%  '
Set-TextValue "E17" '  -3.64%  '
Set-TextValue "D18" '3.193.86'
Set-TextValue "E18" '  -2.90%  '
Set-TextValue "D19" '420.79'
Set-TextValue "E19" '  -1.02%  '
Set-TextValue "D20" '13.01'
Set-TextValue "E20" '  -0.94%  '
Set-TextValue "D21" '5.37'
Set-TextValue "E21" '  -3.13%  '
Set-TextValue "D22" '7.17'
Set-TextValue "E22" '  -2.84%  '
Set-TextValue "E23" '  -0.05%  '
Set-TextValue "D24" '70.33'
Set-TextValue "E24" '  -1.80%  '
Set-TextValue "E25" '  +0.02%  '
Set-TextValue "E26" '  +2.49%  '
Set-TextValue "E27" '  -2.99%  '
Set-TextValue "E28" '  -7.50%  '
Set-TextValue "D29" '8.76'
Set-TextValue "E30" '  -0.50%  '

# --- Rows 31/32 swap places: PancakeSwap <-> EthereumClassic ---
Set-TextValue "B31" 'EthereumClassic'
Set-TextValue "C31" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D31" '21.82'
Set-TextValue "E31" '  -2.49%  '
Set-TextValue "B32" 'PancakeSwap'
Set-TextValue "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '1.83'
Set-TextValue "E32" '  -4.58%  '

Set-TextValue "E33" '  -0.10%  '
Set-TextValue "D34" '5.06'
Set-TextValue "E34" '  -2.18%  '
Set-TextValue "E35" '  -2.84%  '
Set-TextValue "D36" '157.32'
Set-TextValue "E36" '  -2.05%  '
Set-TextValue "E37" '  -4.30%  '
Set-TextValue "D39" '2.713.62'
Set-TextValue "E39" '  -5.24%  '
Set-TextValue "E40" '  -4.83%  '

# --- Rows 41/42 swap places: EnergySwap <-> Filecoin ---
Set-TextValue "B41" 'Filecoin'
Set-TextValue "C41" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D41" '4.24'
Set-TextValue "E41" '  -1.23%  '
Set-TextValue "B42" 'EnergySwap'
Set-TextValue "C42" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D42" '24.37'
Set-TextValue "E42" '  -7.71%  '

Set-TextValue "D43" '39.14'
Set-TextValue "E43" '  -1.63%  '
Set-TextValue "E44" '  -5.43%  '
Set-TextValue "D46" '5.54'
Set-TextValue "E46" '  -6.25%  '
Set-TextValue "E47" '  -2.94%  '
Set-TextValue "E48" '  -6.73%  '
Set-TextValue "D49" '21.44'
Set-TextValue "E49" '  -7.35%  '
Set-TextValue "E50" '  -5.71%  '
Set-TextValue "D51" '0.998'
Set-TextValue "E51" '  -0.19%  '
